# Update refrigerators.xlsx "survey" sheet: add a new filtering "if / assign
# model_row_id / end if" block ahead of the existing model_row_id assignment,
# and rewrite that assignment to use assign('model_row_id', ...) instead of
# returning a bare value (commit message: "Update xlsx again, add filtering").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert 3 blank rows before the old row 16 ("end screen"), pushing everything
# below down by 3 (old row 16 -> new row 19, ..., old row 28 -> new row 31).
$ws.Rows("16:18").Insert()
$ws.Rows("16:18").RowHeight = 16

# New row 16: "if" clause guarding the new assignment block.
$ws.Range("B16").Value = "if"
$ws.Range("C16").Value = 0

# New row 17: a "text" field definition named model_row_id (placeholder row
# added by the form generator).
$ws.Range("D17").Value = "text"
$ws.Range("F17").Value = "model_row_id"

# New row 18: closes the "if" opened on row 16.
$ws.Range("B18").Value = "end if"

# The pre-existing "assign" row (old row 18) is now row 21. Update its name
# and calculation so the JS body explicitly calls assign('model_row_id', x)
# instead of returning the value, renaming the field to sidestep the
# generator's de-duplication ("model_row_id_wtf_xlsxgenerator_stop_that").
$ws.Range("F21").Value = "model_row_id_wtf_xlsxgenerator_stop_that"
$ws.Range("N21").Value = "(function() {var x = 0; if (data('common_models') === 'other') { x = data('model_other'); } else { x = data('common_models');}; assign('model_row_id', x); })()"

# Move the selection roughly where the author's cursor ended up after editing.
$ws.Range("N22").Select()
